{"js": "// Progress report text revisions (Week 6 update), applied as a series of\n// targeted find & replace operations on context.document.body so that the\n// surrounding run formatting (bold names, the red/bold/italic [NOTE]\n// paragraph, etc.) is left untouched.\n\nasync function replaceOnce(body, before, after) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(before) + \" but found \" + results.items.length\n    );\n  }\n  results.items[0].insertText(after, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// Mohand Ferawana paragraph: switch from first-person to third-person,\n// and rewrite the closing sentence about next steps.\nawait replaceOnce(\n  body,\n  \" I worked on the Magnetic strip card reader. This is implemented in the first stage of the main operations procedure for User authentication. Basically, I created a page on the web application where administrative users can swipe student cards from a text field and auto trigger a server request for valid and existing user information before proceeding to the second stage of the main operations procedure which is David\\u2019s item scanning/selection functionality. I made this possible by simply writing a PHP script to extract the data needed from the student card in order to query the database through the web application. After which the returned information is presented in the web page for user verification. I completed this functionality on the 4th of march and no major problems have been encountered so far. This week, I\\u2019ll be meeting with my team again to know the next step for us as a group and for me individually.\",\n  \" Mohand worked on the Magnetic strip card reader. This is implemented in the first stage of the main operations procedure for User authentication. Basically, he created a page on the web application where administrative users can swipe student cards from a text field and auto trigger a server request for valid and existing user information before proceeding to the second stage of the main operations procedure which is David\\u2019s item scanning/selection functionality. He made this possible by extracting the data needed from the student card in order to query the database through the web application. After which the returned information is presented in the web page for user verification. He successfully completed this functionality on the 4th of March and no major problems have been encountered so far. This week, during the group meeting, the team will decide the next step to focus on.\"\n);\n\n// [NOTE] paragraph: \"item rental process\" -> \"item check-out process\", plus\n// a couple of small wording tweaks.\nawait replaceOnce(\n  body,\n  \" is a 3-step procedure we\\u2019re implementing into our web application in order to complete an item rental process. These steps include User Authentication, Item Scanning/Selection and Quantity Selection and Approval.\",\n  \" is a 3-step procedure we\\u2019re implementing into our web application in order to complete the item check-out process. These three steps include User Authentication, Item scanning/selection and Quantity Selection (Summary) and Approval.\"\n);\n\n// Ifeoluwa David paragraph: \"My focus\" -> \"Personally, the focus\".\nawait replaceOnce(body, \"My focus has been on\", \"Personally, the focus has been on\");\n\n// Tosin Ajayi paragraph: switch from first-person to third-person\n// throughout, and rework the final few sentences about Mohand taking over\n// registration.\nawait replaceOnce(body, \" So far, my role in the project \", \" Tosin\\u2019s role in the project \");\n\nawait replaceOnce(\n  body,\n  \"I have successfully completed the login and registration functionality; however, I\\u2019m currently working on implementing a simple validation process for the registration and profile update aspect. This will ensure that the data entered is valid enough to be entered into the database. Based off of Austin\\u2019s email and the agreement with the parts crib, I have also decided that upon completion of this stage. I\\u2019ll be implementing Mohand\\u2019s student card authentication feature into the registration process simply because every card\\u2019s data needs to be extracted and registered at first, before it can be considered valid for user authentication.\",\n  \"He successfully completed the login and registration functionality; however, he\\u2019s currently working on implementing an input validation process for the user registration and profile update aspect. This will ensure that the data entered is valid enough to be entered into the database. Based off of Austin\\u2019s email and the agreement with the parts crib, the group has also decided that upon completion of this stage, Mohand will be taking over the registration process, in order to implement the student card authentication feature using the magnetic stripe card reader. Simply because every card\\u2019s data needs to be extracted and registered at first, before it can be considered valid for user authentication during the item check-out process.\"\n);\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Once($doc, [string]$searchText, [string]$replaceText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($searchText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nReplace-Once $d ' I worked on the Magnetic strip card reader. This is implemented in the first stage of the main operations procedure for User authentication. Basically, I created a page on the web application where administrative users can swipe student cards from a text field and auto trigger a server request for valid and existing user information before proceeding to the second stage of the main operations procedure which is David\u2019s item scanning/selection functionality. I made this possible by simply writing a PHP script to extract the data needed from the student card in order to query the database through the web application. After which the returned information is presented in the web page for user verification. I completed this functionality on the 4th of march and no major problems have been encountered so far. This week, I\u2019ll be meeting with my team again to know the next step for us as a group and for me individually.' ' Mohand worked on the Magnetic strip card reader. This is implemented in the first stage of the main operations procedure for User authentication. Basically, he created a page on the web application where administrative users can swipe student cards from a text field and auto trigger a server request for valid and existing user information before proceeding to the second stage of the main operations procedure which is David\u2019s item scanning/selection functionality. He made this possible by extracting the data needed from the student card in order to query the database through the web application. After which the returned information is presented in the web page for user verification. He successfully completed this functionality on the 4th of March and no major problems have been encountered so far. This week, during the group meeting, the team will decide the next step to focus on.'\nReplace-Once $d ' is a 3-step procedure we\u2019re implementing into our web application in order to complete an item rental process. These steps include User Authentication, Item Scanning/Selection and Quantity Selection and Approval.' ' is a 3-step procedure we\u2019re implementing into our web application in order to complete the item check-out process. These three steps include User Authentication, Item scanning/selection and Quantity Selection (Summary) and Approval.'\nReplace-Once $d 'My focus has been on' 'Personally, the focus has been on'\nReplace-Once $d ' So far, my role in the project ' ' Tosin\u2019s role in the project '\nReplace-Once $d 'I have successfully completed the login and registration functionality; however, I\u2019m currently working on implementing a simple validation process for the registration and profile update aspect. This will ensure that the data entered is valid enough to be entered into the database. Based off of Austin\u2019s email and the agreement with the parts crib, I have also decided that upon completion of this stage. I\u2019ll be implementing Mohand\u2019s student card authentication feature into the registration process simply because every card\u2019s data needs to be extracted and registered at first, before it can be considered valid for user authentication.' 'He successfully completed the login and registration functionality; however, he\u2019s currently working on implementing an input validation process for the user registration and profile update aspect. This will ensure that the data entered is valid enough to be entered into the database. Based off of Austin\u2019s email and the agreement with the parts crib, the group has also decided that upon completion of this stage, Mohand will be taking over the registration process, in order to implement the student card authentication feature using the magnetic stripe card reader. Simply because every card\u2019s data needs to be extracted and registered at first, before it can be considered valid for user authentication during the item check-out process.'\n"}
